$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: add new field row (row 81) describing the new "Week type" field ---
$ws1.Range("A81").Value = 89
$ws1.Range("B81").Value = "Week type"
$ws1.Range("C81").Value = "weektype"
$ws1.Range("D81").Value = "varchar"
$ws1.Range("E81").Value = "weeks"

# --- Sheet2: insert a row for the new field into the "Weeks" summary block ---
$ws2.Rows.Item(51).Insert()

# Copy the border/style formatting that used to sit on the last row of the
# block (row 50, E:F) onto the freshly inserted row so the "closing" border
# moves down together with the new last data row.
$ws2.Range("E50:F50").Copy($ws2.Range("E51:F51"))

# Fill in the new row's content (mirrors the other rows in this block, which
# pull their field name / description from Sheet1 via formulas).
$ws2.Range("E51").Formula = "=Sheet1!C81"
$ws2.Range("F51").Formula = "=Sheet1!B81"

# Row 50 (end/End Date) is no longer the last row of the block, so it takes
# the "continuing" style instead (matches row 49's style/border).
$ws2.Range("E49:F49").Copy($ws2.Range("E50:F50"))
$ws2.Range("E50").Formula = "=Sheet1!C80"
$ws2.Range("F50").Formula = "=Sheet1!B80"

# --- Restore view/selection state ---
$ws1.Range("A82").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("E54").Select() | Out-Null

Write-Output "done"
